# Updated symbol list on Mon Dec 19 02:27:16 UTC 2022 with GitHub Actions
#
# Refreshes the cryptocurrency price snapshot in column D and bumps the
# "Hora" (hour) marker in column G from "1" to "2" for every data row
# (rows 2-51). Values are written with NumberFormat "@" (Text) first so
# Excel keeps them as literal strings (preserving exact decimal text such
# as trailing/leading zeros) instead of silently re-parsing them as
# numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Price" (column D) values, keyed by row number. Rows not listed here
# (24, 27-41, 51) keep their existing Price value unchanged.
$priceUpdates = @{
    2  = "251.08"
    3  = "22.06"
    4  = "5.563"
    5  = "0.05676"
    6  = "6.486"
    7  = "0.8042"
    8  = "1.051"
    9  = "0.1435"
    10 = "0.07308"
    11 = "0.03197"
    12 = "0.02944"
    13 = "0.09262"
    14 = "0.001666"
    15 = "3.243"
    16 = "0.04741"
    17 = "0.0005818"
    18 = "0.006430"
    19 = "0.005062"
    20 = "0.001051"
    21 = "0.0001502"
    22 = "0.0003204"
    23 = "4.065"
    25 = "2.114"
    26 = "0.3272"
    42 = "0.002973"
    43 = "0.006908"
    44 = "0.008529"
    45 = "0.00005648"
    46 = "0.00000000751"
    47 = "0.7862"
    48 = "0.01745"
    49 = "0.00002103"
    50 = "0.01011"
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
}

# "Hora" (column G) goes from "1" to "2" on every data row (2 through 51).
for ($row = 2; $row -le 51; $row++) {
    $cell = $ws.Range("G$row")
    $cell.NumberFormat = "@"
    $cell.Value = "2"
}
